$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 28.72393179789166
$ws.Range("C2").Value = 9.74716779482066
$ws.Range("D2").Value = 4.276585547397774
$ws.Range("E2").Value = 9.782116337054427
$ws.Range("F2").Value = 66.80697039126721
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 10.43950258765039
$ws.Range("L2").Value = 12.05500587233369

$ws.Range("B3").Value = 28.61399692459371
$ws.Range("C3").Value = 9.534477755564742
$ws.Range("D3").Value = 4.146541706473907
$ws.Range("E3").Value = 9.794217008429589
$ws.Range("F3").Value = 66.08808740486724
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 10.44178200816691
$ws.Range("L3").Value = 12.09639430508557

$ws.Range("B4").Value = 28.55689149624353
$ws.Range("C4").Value = 9.407063702295341
$ws.Range("D4").Value = 4.064582800198115
$ws.Range("E4").Value = 9.802104641717801
$ws.Range("F4").Value = 65.65189230553317
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 10.44364368826178
$ws.Range("L4").Value = 12.12426386912827

$ws.Range("B5").Value = 28.53625169539993
$ws.Range("C5").Value = 9.35603993239269
$ws.Range("D5").Value = 4.030685727209747
$ws.Range("E5").Value = 9.805434356480255
$ws.Range("F5").Value = 65.47556705947595
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 10.44451826866739
$ws.Range("L5").Value = 12.1362385036599

$ws.Range("B6").Value = 28.53298380999409
$ws.Range("C6").Value = 9.347624623434626
$ws.Range("D6").Value = 4.025028065369135
$ws.Range("E6").Value = 9.805994235110534
$ws.Range("F6").Value = 65.44637804293779
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 10.44467048606404
$ws.Range("L6").Value = 12.13826417438762

$ws.Range("B7").Value = 28.55660246749541
$ws.Range("C7").Value = 9.406371809385623
$ws.Range("D7").Value = 4.064127624221029
$ws.Range("E7").Value = 9.802149079559392
$ws.Range("F7").Value = 65.64950838521656
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 10.44365501411513
$ws.Range("L7").Value = 12.12442286306488

$ws.Range("B8").Value = 28.68388029844533
$ws.Range("C8").Value = 9.673232263809201
$ws.Range("D8").Value = 4.23220353672086
$ws.Range("E8").Value = 9.786193835926969
$ws.Range("F8").Value = 66.55809279230851
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 10.44019245891559
$ws.Range("L8").Value = 12.0687665288739

$ws.Range("B9").Value = 29.0150082774443
$ws.Range("C9").Value = 10.21740486231168
$ws.Range("D9").Value = 4.543739281661184
$ws.Range("E9").Value = 9.758523066601363
$ws.Range("F9").Value = 68.37569425621032
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 10.43708166561301
$ws.Range("L9").Value = 11.97913038574781

$ws.Range("B10").Value = 29.30647077796392
$ws.Range("C10").Value = 10.62419420524044
$ws.Range("D10").Value = 4.760050035971784
$ws.Range("E10").Value = 9.740378253038157
$ws.Range("F10").Value = 69.72556874423755
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 10.43705738439145
$ws.Range("L10").Value = 11.92518198382256

$ws.Range("B11").Value = 29.44912663668646
$ws.Range("C11").Value = 10.80961680102561
$ws.Range("D11").Value = 4.855444979789064
$ws.Range("E11").Value = 9.732593851312
$ws.Range("F11").Value = 70.34120785380037
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 10.43754127291909
$ws.Range("L11").Value = 11.90322837060306

$ws.Range("B12").Value = 29.50455575985173
$ws.Range("C12").Value = 10.87979501633352
$ws.Range("D12").Value = 4.891116087259291
$ws.Range("E12").Value = 9.729713325611939
$ws.Range("F12").Value = 70.57442545900111
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 10.43779597958907
$ws.Range("L12").Value = 11.89528750099568

$ws.Range("B13").Value = 29.49255612304256
$ws.Range("C13").Value = 10.86468393055557
$ws.Range("D13").Value = 4.883454123037623
$ws.Range("E13").Value = 9.73033071166264
$ws.Range("F13").Value = 70.5241961283603
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 10.43773794065436
$ws.Range("L13").Value = 11.89698113717173

$ws.Range("B14").Value = 29.45365883250085
$ws.Range("C14").Value = 10.81539165394378
$ws.Range("D14").Value = 4.858388849998764
$ws.Range("E14").Value = 9.732355522438702
$ws.Range("F14").Value = 70.36039382115879
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 10.43756079375465
$ws.Range("L14").Value = 11.90256760323457

$ws.Range("B15").Value = 29.43001532191356
$ws.Range("C15").Value = 10.78519124852999
$ws.Range("D15").Value = 4.84297609327333
$ws.Range("E15").Value = 9.733604527880974
$ws.Range("F15").Value = 70.26006761568047
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 10.43746160228405
$ws.Range("L15").Value = 11.90603799447779

$ws.Range("B16").Value = 29.29734718437883
$ws.Range("C16").Value = 10.61207600965625
$ws.Range("D16").Value = 4.75375344894919
$ws.Range("E16").Value = 9.740896409706489
$ws.Range("F16").Value = 69.68535515067938
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 10.43703574597426
$ws.Range("L16").Value = 11.92666880861474

$ws.Range("B17").Value = 29.21851124231646
$ws.Range("C17").Value = 10.50590448819024
$ws.Range("D17").Value = 4.698232947265764
$ws.Range("E17").Value = 9.745489851148138
$ws.Range("F17").Value = 69.33309253678166
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 10.43690147392813
$ws.Range("L17").Value = 11.93998820417358

$ws.Range("B18").Value = 29.17411710747684
$ws.Range("C18").Value = 10.44487758185211
$ws.Range("D18").Value = 4.666017177845855
$ws.Range("E18").Value = 9.748176110979111
$ws.Range("F18").Value = 69.13063682386466
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 10.43687081210848
$ws.Range("L18").Value = 11.94789272183021

$ws.Range("B19").Value = 29.15925036506711
$ws.Range("C19").Value = 10.42422466585992
$ws.Range("D19").Value = 4.655061731176864
$ws.Range("E19").Value = 9.749093238601347
$ws.Range("F19").Value = 69.06211988003183
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 10.43686841952405
$ws.Range("L19").Value = 11.95061087480817

$ws.Range("B20").Value = 29.22680539483806
$ws.Range("C20").Value = 10.51720314736451
$ws.Range("D20").Value = 4.704172526115119
$ws.Range("E20").Value = 9.744996295572056
$ws.Range("F20").Value = 69.37057630592606
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 10.43691094565244
$ws.Range("L20").Value = 11.93854512221052

$ws.Range("B21").Value = 29.46504601294554
$ws.Range("C21").Value = 10.8298717002734
$ws.Range("D21").Value = 4.865763572822805
$ws.Range("E21").Value = 9.731758963379548
$ws.Range("F21").Value = 70.40850522423531
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 10.43761088409564
$ws.Range("L21").Value = 11.9009166106804

$ws.Range("B22").Value = 29.62893925813421
$ws.Range("C22").Value = 11.0339629591102
$ws.Range("D22").Value = 4.968724881732752
$ws.Range("E22").Value = 9.723499501008849
$ws.Range("F22").Value = 71.08730405201163
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 10.43848504465961
$ws.Range("L22").Value = 11.87849550388937

$ws.Range("B23").Value = 29.54073077834026
$ws.Range("C23").Value = 10.92508763105649
$ws.Range("D23").Value = 4.91402100698331
$ws.Range("E23").Value = 9.72787196859813
$ws.Range("F23").Value = 70.72502026566926
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 10.43798026187737
$ws.Range("L23").Value = 11.89026327884853

$ws.Range("B24").Value = 29.2230527097328
$ws.Range("C24").Value = 10.51209497761958
$ws.Range("D24").Value = 4.701488163814618
$ws.Range("E24").Value = 9.745219290524117
$ws.Range("F24").Value = 69.35362968172275
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 10.43690651856879
$ws.Range("L24").Value = 11.9391967701826

$ws.Range("B25").Value = 28.91684298281355
$ws.Range("C25").Value = 10.06856034083797
$ws.Range("D25").Value = 4.461565550641651
$ws.Range("E25").Value = 9.765623593518544
$ws.Range("F25").Value = 67.88089990786584
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 10.43752743128697
$ws.Range("L25").Value = 12.00128941426054
